# Updates the IFRS consolidated financial-summary rows (2014/12 .. 2021/12(E))
# for 한화손해보험: the figures were re-pulled (values now correct per the fixed
# ingestion pipeline), several no-longer-reported breakdown columns (N/A "비지배"
# and "FCF" style columns) are cleared out entirely for the affected rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (annual period column)
$ws.Cells.Item(2, 4).Value = 53292   # D2
$ws.Cells.Item(2, 5).Value = 277   # E2
$ws.Cells.Item(2, 6).Value = 277   # F2
$ws.Cells.Item(2, 7).Value = 268   # G2
$ws.Cells.Item(2, 8).Value = 129   # H2
$ws.Cells.Item(2, 9).Value = 129   # I2
$ws.Cells.Item(2, 10).Value = 0   # J2
$ws.Cells.Item(2, 11).Value = 103316   # K2
$ws.Cells.Item(2, 12).Value = 97366   # L2
$ws.Cells.Item(2, 13).Value = 5950   # M2
$ws.Cells.Item(2, 14).Value = 5950   # N2
$ws.Cells.Item(2, 15).ClearContents()   # O2
$ws.Cells.Item(2, 16).Value = 4537   # P2
$ws.Cells.Item(2, 17).Value = 9176   # Q2
$ws.Cells.Item(2, 18).Value = -8463   # R2
$ws.Cells.Item(2, 19).Value = 0   # S2
$ws.Cells.Item(2, 20).Value = 86   # T2
$ws.Cells.Item(2, 21).ClearContents()   # U2
$ws.Cells.Item(2, 22).Value = 1944   # V2
$ws.Cells.Item(2, 23).Value = 0.52   # W2
$ws.Cells.Item(2, 24).Value = 0.24   # X2
$ws.Cells.Item(2, 25).Value = 2.25   # Y2
$ws.Cells.Item(2, 26).Value = 0.13   # Z2
$ws.Cells.Item(2, 27).Value = 1636.29   # AA2
$ws.Cells.Item(2, 28).Value = 31.16   # AB2
$ws.Cells.Item(2, 29).Value = 140   # AC2
$ws.Cells.Item(2, 30).Value = 32.03   # AD2
$ws.Cells.Item(2, 31).Value = 6498   # AE2
$ws.Cells.Item(2, 32).Value = 0.69   # AF2
$ws.Cells.Item(2, 33).Value = 0   # AG2
$ws.Cells.Item(2, 34).Value = 0   # AH2
$ws.Cells.Item(2, 35).Value = 0   # AI2
$ws.Cells.Item(2, 36).Value = 91579218   # AJ2

# Row 3 (annual period column)
$ws.Cells.Item(3, 4).Value = 57471   # D3
$ws.Cells.Item(3, 5).Value = 1191   # E3
$ws.Cells.Item(3, 6).Value = 1191   # F3
$ws.Cells.Item(3, 7).Value = 1171   # G3
$ws.Cells.Item(3, 8).Value = 958   # H3
$ws.Cells.Item(3, 9).Value = 958   # I3
$ws.Cells.Item(3, 10).ClearContents()   # J3
$ws.Cells.Item(3, 11).Value = 118034   # K3
$ws.Cells.Item(3, 12).Value = 109425   # L3
$ws.Cells.Item(3, 13).Value = 8609   # M3
$ws.Cells.Item(3, 14).Value = 8609   # N3
$ws.Cells.Item(3, 15).ClearContents()   # O3
$ws.Cells.Item(3, 16).Value = 4537   # P3
$ws.Cells.Item(3, 17).Value = 3316   # Q3
$ws.Cells.Item(3, 18).Value = -3442   # R3
$ws.Cells.Item(3, 19).Value = -300   # S3
$ws.Cells.Item(3, 20).Value = 93   # T3
$ws.Cells.Item(3, 21).ClearContents()   # U3
$ws.Cells.Item(3, 22).Value = 1646   # V3
$ws.Cells.Item(3, 23).Value = 2.07   # W3
$ws.Cells.Item(3, 24).Value = 1.67   # X3
$ws.Cells.Item(3, 25).Value = 13.15   # Y3
$ws.Cells.Item(3, 26).Value = 0.87   # Z3
$ws.Cells.Item(3, 27).Value = 1271.04   # AA3
$ws.Cells.Item(3, 28).Value = 89.75   # AB3
$ws.Cells.Item(3, 29).Value = 1046   # AC3
$ws.Cells.Item(3, 30).Value = 7.3   # AD3
$ws.Cells.Item(3, 31).Value = 9401   # AE3
$ws.Cells.Item(3, 32).Value = 0.81   # AF3
$ws.Cells.Item(3, 33).Value = 69   # AG3
$ws.Cells.Item(3, 34).Value = 0.91   # AH3
$ws.Cells.Item(3, 35).Value = 6.63   # AI3
$ws.Cells.Item(3, 36).Value = 91579218   # AJ3

# Row 4 (annual period column)
$ws.Cells.Item(4, 4).Value = 61931   # D4
$ws.Cells.Item(4, 5).Value = 1354   # E4
$ws.Cells.Item(4, 6).Value = 1354   # F4
$ws.Cells.Item(4, 7).Value = 1341   # G4
$ws.Cells.Item(4, 8).Value = 1116   # H4
$ws.Cells.Item(4, 9).Value = 1116   # I4
$ws.Cells.Item(4, 10).ClearContents()   # J4
$ws.Cells.Item(4, 11).Value = 132833   # K4
$ws.Cells.Item(4, 12).Value = 123989   # L4
$ws.Cells.Item(4, 13).Value = 8845   # M4
$ws.Cells.Item(4, 14).Value = 8845   # N4
$ws.Cells.Item(4, 15).ClearContents()   # O4
$ws.Cells.Item(4, 16).Value = 4537   # P4
$ws.Cells.Item(4, 17).Value = 5161   # Q4
$ws.Cells.Item(4, 18).Value = -5748   # R4
$ws.Cells.Item(4, 19).Value = 811   # S4
$ws.Cells.Item(4, 20).Value = 1734   # T4
$ws.Cells.Item(4, 21).ClearContents()   # U4
$ws.Cells.Item(4, 22).Value = 2523   # V4
$ws.Cells.Item(4, 23).Value = 2.19   # W4
$ws.Cells.Item(4, 24).Value = 1.8   # X4
$ws.Cells.Item(4, 25).Value = 12.79   # Y4
$ws.Cells.Item(4, 26).Value = 0.89   # Z4
$ws.Cells.Item(4, 27).Value = 1401.81   # AA4
$ws.Cells.Item(4, 28).Value = 94.95   # AB4
$ws.Cells.Item(4, 29).Value = 1219   # AC4
$ws.Cells.Item(4, 30).Value = 5.81   # AD4
$ws.Cells.Item(4, 31).Value = 9658   # AE4
$ws.Cells.Item(4, 32).Value = 0.73   # AF4
$ws.Cells.Item(4, 33).Value = 99   # AG4
$ws.Cells.Item(4, 34).Value = 1.4   # AH4
$ws.Cells.Item(4, 35).Value = 8.13   # AI4
$ws.Cells.Item(4, 36).Value = 91579218   # AJ4

# Row 5 (annual period column)
$ws.Cells.Item(5, 4).Value = 69318   # D5
$ws.Cells.Item(5, 5).Value = 1975   # E5
$ws.Cells.Item(5, 6).Value = 1975   # F5
$ws.Cells.Item(5, 7).Value = 1939   # G5
$ws.Cells.Item(5, 8).Value = 1476   # H5
$ws.Cells.Item(5, 9).Value = 1476   # I5
$ws.Cells.Item(5, 10).ClearContents()   # J5
$ws.Cells.Item(5, 11).Value = 148906   # K5
$ws.Cells.Item(5, 12).Value = 137062   # L5
$ws.Cells.Item(5, 13).Value = 11843   # M5
$ws.Cells.Item(5, 14).Value = 11843   # N5
$ws.Cells.Item(5, 15).ClearContents()   # O5
$ws.Cells.Item(5, 16).Value = 5837   # P5
$ws.Cells.Item(5, 17).Value = 9386   # Q5
$ws.Cells.Item(5, 18).Value = -10642   # R5
$ws.Cells.Item(5, 19).Value = 1752   # S5
$ws.Cells.Item(5, 20).Value = 123   # T5
$ws.Cells.Item(5, 21).ClearContents()   # U5
$ws.Cells.Item(5, 22).Value = 2174   # V5
$ws.Cells.Item(5, 23).Value = 2.85   # W5
$ws.Cells.Item(5, 24).Value = 2.13   # X5
$ws.Cells.Item(5, 25).Value = 14.27   # Y5
$ws.Cells.Item(5, 26).Value = 1.05   # Z5
$ws.Cells.Item(5, 27).Value = 1157.32   # AA5
$ws.Cells.Item(5, 28).Value = 102.9   # AB5
$ws.Cells.Item(5, 29).Value = 1551   # AC5
$ws.Cells.Item(5, 30).Value = 5.25   # AD5
$ws.Cells.Item(5, 31).Value = 10145   # AE5
$ws.Cells.Item(5, 32).Value = 0.8   # AF5
$ws.Cells.Item(5, 33).Value = 150   # AG5
$ws.Cells.Item(5, 34).Value = 1.84   # AH5
$ws.Cells.Item(5, 35).Value = 11.86   # AI5
$ws.Cells.Item(5, 36).Value = 116738915   # AJ5

# Row 6 (annual period column)
$ws.Cells.Item(6, 4).Value = 74238   # D6
$ws.Cells.Item(6, 5).Value = 1105   # E6
$ws.Cells.Item(6, 6).Value = 1105   # F6
$ws.Cells.Item(6, 7).Value = 1131   # G6
$ws.Cells.Item(6, 8).Value = 818   # H6
$ws.Cells.Item(6, 9).Value = 818   # I6
$ws.Cells.Item(6, 10).ClearContents()   # J6
$ws.Cells.Item(6, 11).Value = 167339   # K6
$ws.Cells.Item(6, 12).Value = 153183   # L6
$ws.Cells.Item(6, 13).Value = 14156   # M6
$ws.Cells.Item(6, 14).Value = 14156   # N6
$ws.Cells.Item(6, 15).ClearContents()   # O6
$ws.Cells.Item(6, 16).Value = 5837   # P6
$ws.Cells.Item(6, 17).Value = 7107   # Q6
$ws.Cells.Item(6, 18).Value = -12223   # R6
$ws.Cells.Item(6, 19).Value = 5156   # S6
$ws.Cells.Item(6, 20).Value = 313   # T6
$ws.Cells.Item(6, 21).ClearContents()   # U6
$ws.Cells.Item(6, 22).Value = 5664   # V6
$ws.Cells.Item(6, 23).Value = 1.49   # W6
$ws.Cells.Item(6, 24).Value = 1.1   # X6
$ws.Cells.Item(6, 25).Value = 6.29   # Y6
$ws.Cells.Item(6, 26).Value = 0.52   # Z6
$ws.Cells.Item(6, 27).Value = 1082.14   # AA6
$ws.Cells.Item(6, 28).Value = 142.52   # AB6
$ws.Cells.Item(6, 29).Value = 700   # AC6
$ws.Cells.Item(6, 30).Value = 8.44   # AD6
$ws.Cells.Item(6, 31).Value = 12126   # AE6
$ws.Cells.Item(6, 32).Value = 0.49   # AF6
$ws.Cells.Item(6, 33).Value = 130   # AG6
$ws.Cells.Item(6, 34).Value = 2.2   # AH6
$ws.Cells.Item(6, 35).Value = 18.56   # AI6
$ws.Cells.Item(6, 36).Value = 116738915   # AJ6

# Row 7 (estimate period column)
$ws.Cells.Item(7, 4).Value = 45660   # D7
$ws.Cells.Item(7, 5).ClearContents()   # E7
$ws.Cells.Item(7, 6).ClearContents()   # F7
$ws.Cells.Item(7, 7).Value = -380   # G7
$ws.Cells.Item(7, 8).Value = -240   # H7
$ws.Cells.Item(7, 9).Value = -240   # I7
$ws.Cells.Item(7, 10).ClearContents()   # J7
$ws.Cells.Item(7, 11).Value = 182460   # K7
$ws.Cells.Item(7, 12).Value = 167040   # L7
$ws.Cells.Item(7, 13).Value = 15420   # M7
$ws.Cells.Item(7, 14).Value = 15130   # N7
$ws.Cells.Item(7, 15).ClearContents()   # O7
$ws.Cells.Item(7, 16).ClearContents()   # P7
$ws.Cells.Item(7, 17).ClearContents()   # Q7
$ws.Cells.Item(7, 18).ClearContents()   # R7
$ws.Cells.Item(7, 19).ClearContents()   # S7
$ws.Cells.Item(7, 20).ClearContents()   # T7
$ws.Cells.Item(7, 21).ClearContents()   # U7
$ws.Cells.Item(7, 22).ClearContents()   # V7
$ws.Cells.Item(7, 23).ClearContents()   # W7
$ws.Cells.Item(7, 24).Value = -0.53   # X7
$ws.Cells.Item(7, 25).Value = -1.64   # Y7
$ws.Cells.Item(7, 26).Value = -0.14   # Z7
$ws.Cells.Item(7, 27).Value = 1083.27   # AA7
$ws.Cells.Item(7, 28).ClearContents()   # AB7
$ws.Cells.Item(7, 29).Value = -206   # AC7
$ws.Cells.Item(7, 30).Value = -12.4   # AD7
$ws.Cells.Item(7, 31).Value = 12961   # AE7
$ws.Cells.Item(7, 32).Value = 0.2   # AF7
$ws.Cells.Item(7, 33).Value = 20   # AG7
$ws.Cells.Item(7, 34).Value = 0.78   # AH7
$ws.Cells.Item(7, 35).Value = -9.73   # AI7
$ws.Cells.Item(7, 36).ClearContents()   # AJ7

# Row 8 (estimate period column)
$ws.Cells.Item(8, 4).Value = 47610   # D8
$ws.Cells.Item(8, 5).ClearContents()   # E8
$ws.Cells.Item(8, 6).ClearContents()   # F8
$ws.Cells.Item(8, 7).Value = 570   # G8
$ws.Cells.Item(8, 8).Value = 440   # H8
$ws.Cells.Item(8, 9).Value = 440   # I8
$ws.Cells.Item(8, 10).ClearContents()   # J8
$ws.Cells.Item(8, 11).Value = 195050   # K8
$ws.Cells.Item(8, 12).Value = 179200   # L8
$ws.Cells.Item(8, 13).Value = 15850   # M8
$ws.Cells.Item(8, 14).Value = 15560   # N8
$ws.Cells.Item(8, 15).ClearContents()   # O8
$ws.Cells.Item(8, 16).ClearContents()   # P8
$ws.Cells.Item(8, 17).ClearContents()   # Q8
$ws.Cells.Item(8, 18).ClearContents()   # R8
$ws.Cells.Item(8, 19).ClearContents()   # S8
$ws.Cells.Item(8, 20).ClearContents()   # T8
$ws.Cells.Item(8, 21).ClearContents()   # U8
$ws.Cells.Item(8, 22).ClearContents()   # V8
$ws.Cells.Item(8, 23).ClearContents()   # W8
$ws.Cells.Item(8, 24).Value = 0.92   # X8
$ws.Cells.Item(8, 25).Value = 2.87   # Y8
$ws.Cells.Item(8, 26).Value = 0.23   # Z8
$ws.Cells.Item(8, 27).Value = 1130.6   # AA8
$ws.Cells.Item(8, 28).ClearContents()   # AB8
$ws.Cells.Item(8, 29).Value = 377   # AC8
$ws.Cells.Item(8, 30).Value = 6.21   # AD8
$ws.Cells.Item(8, 31).Value = 13329   # AE8
$ws.Cells.Item(8, 32).Value = 0.18   # AF8
$ws.Cells.Item(8, 33).Value = 44   # AG8
$ws.Cells.Item(8, 34).Value = 1.88   # AH8
$ws.Cells.Item(8, 35).Value = 11.67   # AI8
$ws.Cells.Item(8, 36).ClearContents()   # AJ8

# Row 9 (estimate period column)
$ws.Cells.Item(9, 4).Value = 48960   # D9
$ws.Cells.Item(9, 5).ClearContents()   # E9
$ws.Cells.Item(9, 6).ClearContents()   # F9
$ws.Cells.Item(9, 7).Value = 820   # G9
$ws.Cells.Item(9, 8).Value = 620   # H9
$ws.Cells.Item(9, 9).Value = 620   # I9
$ws.Cells.Item(9, 10).ClearContents()   # J9
$ws.Cells.Item(9, 11).Value = 203440   # K9
$ws.Cells.Item(9, 12).Value = 187040   # L9
$ws.Cells.Item(9, 13).Value = 16400   # M9
$ws.Cells.Item(9, 14).Value = 16100   # N9
$ws.Cells.Item(9, 15).ClearContents()   # O9
$ws.Cells.Item(9, 16).ClearContents()   # P9
$ws.Cells.Item(9, 17).ClearContents()   # Q9
$ws.Cells.Item(9, 18).ClearContents()   # R9
$ws.Cells.Item(9, 19).ClearContents()   # S9
$ws.Cells.Item(9, 20).ClearContents()   # T9
$ws.Cells.Item(9, 21).ClearContents()   # U9
$ws.Cells.Item(9, 22).ClearContents()   # V9
$ws.Cells.Item(9, 23).ClearContents()   # W9
$ws.Cells.Item(9, 24).Value = 1.27   # X9
$ws.Cells.Item(9, 25).Value = 3.92   # Y9
$ws.Cells.Item(9, 26).Value = 0.31   # Z9
$ws.Cells.Item(9, 27).Value = 1140.49   # AA9
$ws.Cells.Item(9, 28).ClearContents()   # AB9
$ws.Cells.Item(9, 29).Value = 531   # AC9
$ws.Cells.Item(9, 30).Value = 4.41   # AD9
$ws.Cells.Item(9, 31).Value = 13791   # AE9
$ws.Cells.Item(9, 32).Value = 0.17   # AF9
$ws.Cells.Item(9, 33).Value = 39   # AG9
$ws.Cells.Item(9, 34).Value = 1.67   # AH9
$ws.Cells.Item(9, 35).Value = 7.34   # AI9
$ws.Cells.Item(9, 36).ClearContents()   # AJ9
